$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D width change (22 -> 21)
$ws.Columns.Item(4).ColumnWidth = 20.14

# Row 3
$ws.Range("B3").Value = -14
$ws.Range("C3").Value = 5.1
$ws.Range("D3").Value = "Cloudy "

# Row 4
$ws.Range("B4").Value = -14.4
$ws.Range("C4").Value = 6.3
$ws.Range("D4").Value = "Light freezing rain"

# Row 5
$ws.Range("B5").Value = -5.2
$ws.Range("C5").Value = 6.8
$ws.Range("D5").Value = "Heavy snow"

# Row 6
$ws.Range("B6").Value = -10.4
$ws.Range("C6").Value = 4.4
$ws.Range("D6").Value = "Light snow"

# Row 7
$ws.Range("B7").Value = -23.6
$ws.Range("C7").Value = 2.4
$ws.Range("D7").Value = "Freezing fog"

# Row 8
$ws.Range("B8").Value = -19
$ws.Range("C8").Value = 7
$ws.Range("D8").Value = "Light freezing rain"

# Row 9
$ws.Range("B9").Value = -5.2
$ws.Range("C9").Value = 6
$ws.Range("D9").Value = "Moderate snow"
